$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "67.727.44"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.326.18"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "3.322.03"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "705.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "3.871.82"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "67.757.75"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "3.324.93"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +6.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.97%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "3.701.09"
$ws.Range("E37").Value = "  -5.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "0.0₃0674"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.48%  "
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("E50").Value = "  -5.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.77%  "
